# Apply update: "Förändrad" (C) date bumped from 45671 to 45672 for rows 2-36,
# and rows 35/36 swap their "Beteckning" (A) and "Area (ha)" (G) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad date) for rows 2 through 36: 45671 -> 45672
for ($r = 2; $r -le 36; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45671) {
        $cell.Value2 = 45672
    }
}

# Swap the Beteckning (A) values between rows 35 and 36
$a35 = $ws.Cells.Item(35, 1).Value2
$a36 = $ws.Cells.Item(36, 1).Value2
$ws.Cells.Item(35, 1).Value2 = $a36
$ws.Cells.Item(36, 1).Value2 = $a35

# Swap the Area (ha) (G) values between rows 35 and 36
$g35 = $ws.Cells.Item(35, 7).Value2
$g36 = $ws.Cells.Item(36, 7).Value2
$ws.Cells.Item(35, 7).Value2 = $g36
$ws.Cells.Item(36, 7).Value2 = $g35
